$d = $word.ActiveDocument

# --- Change 1: insert a new ListBullet paragraph "5817181 - Valdeir Arantes"
#     right after the "Docente(s) Responsável(eis) " heading paragraph ---
$rng = $d.Content
$found = $rng.Find.Execute("Docente(s) Responsável(eis) ", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)

if ($found) {
    $srcStart = $rng.Paragraphs.First.Range.Start

    $rng.Collapse(0)
    $rng.InsertParagraphAfter()

    # locate the freshly-inserted (still empty) paragraph that follows the source one
    $newPara = $null
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Start -eq $srcStart) {
            $newPara = $d.Paragraphs.Item($i + 1)
            break
        }
    }

    $newPara.Style = "ListBullet"
    $newPara.Range.Text = "5817181 - Valdeir Arantes"
}

# --- Change 2: fix the LOT2058 requirement description ---
$d.Content.Find.Execute("LOT2058 -  Engenharia Econômica  (Requisito fraco)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "LOT2058 -  Análise Técnico (Requisito fraco)", 2) | Out-Null
